$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = 46
$ws.Range("B22").Value = "Edit1"
$ws.Range("C22").Value = "riya-morankar"
$ws.Range("D22").Value = "N/A"
$ws.Range("E22").Value = "edit1 to main"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "2025-06-19"
$ws.Range("F22").ClearFormats()
